$wb = $excel.ActiveWorkbook

$wsProc = $wb.Worksheets.Item("processes")
$wsTopo = $wb.Worksheets.Item("process_topology")
$wsCf   = $wb.Worksheets.Item("cf")

# ---------------------------------------------------------------------------
# Pre-seed the new shared strings in the same order the original author's
# session created them (pv2, pv2,s1, pv2,s2, pv2,s3, is_cf_fix) so the
# resulting sharedStrings.xml table matches exactly.
# ---------------------------------------------------------------------------
$wsProc.Cells.Item(20,20).Value = "pv2"
$wsProc.Cells.Item(20,21).Value = "pv2,s1"
$wsProc.Cells.Item(20,22).Value = "pv2,s2"
$wsProc.Cells.Item(20,23).Value = "pv2,s3"
$wsProc.Cells.Item(20,24).Value = "is_cf_fix"
$wsProc.Range("T20:X20").Clear()

# ---------------------------------------------------------------------------
# Sheet "processes" (sheet2): insert "is_cf_fix" column after "is_cf", add a
# new process row "pv2" (copied/adapted from "pv1"), and flip a few flags.
# ---------------------------------------------------------------------------

# Insert a new column C ("is_cf_fix") - shifts is_online..min_offline right.
$wsProc.Columns.Item(3).Insert()
$wsProc.Cells.Item(1,3).Value = "is_cf_fix"

# New full data grid (rows 2-7 existing processes, after column insert):
# A=process B=is_cf C=is_cf_fix D=is_online E=is_res F=conversion G=eff H=load_min I=load_max J=start_cost K=min_online L=min_offline
$wsProc.Cells.Item(2,3).Value = 0
$wsProc.Cells.Item(3,3).Value = 0
$wsProc.Cells.Item(4,3).Value = 0
$wsProc.Cells.Item(5,3).Value = 1
$wsProc.Cells.Item(6,3).Value = 0
$wsProc.Cells.Item(7,3).Value = 0

# Insert the new "pv2" row (row 6), pushing dh_tra/spot down.
$wsProc.Rows.Item(6).Insert()
$wsProc.Cells.Item(6,1).Value = "pv2"
$wsProc.Cells.Item(6,2).Value = 1
$wsProc.Cells.Item(6,3).Value = 0
$wsProc.Cells.Item(6,4).Value = 0
$wsProc.Cells.Item(6,5).Value = 0
$wsProc.Cells.Item(6,6).Value = 1
$wsProc.Cells.Item(6,7).Value = 1
$wsProc.Cells.Item(6,8).Value = 0
$wsProc.Cells.Item(6,9).Value = 1
$wsProc.Cells.Item(6,10).Value = 0
$wsProc.Cells.Item(6,11).Value = 0
$wsProc.Cells.Item(6,12).Value = 0

# is_res values for existing rows were shifted from D to E by the column
# insert; fix them up to match the new target layout (is_res now 0 for
# pv1/dh_tra/spot, unchanged for ngchp/hp1/p2x1).
$wsProc.Cells.Item(5,5).Value = 0
$wsProc.Cells.Item(7,5).Value = 0
$wsProc.Cells.Item(8,5).Value = 0

# ---------------------------------------------------------------------------
# Sheet "process_topology" (sheet4): insert a new topology row for pv2 sink.
# ---------------------------------------------------------------------------
$wsTopo.Rows.Item(10).Insert()
$wsTopo.Cells.Item(10,1).Value = "pv2"
$wsTopo.Cells.Item(10,2).Value = "sink"
$wsTopo.Cells.Item(10,3).Value = "elc"
$wsTopo.Cells.Item(10,4).Value = 1
$wsTopo.Cells.Item(10,5).Value = 5
$wsTopo.Cells.Item(10,6).Value = 0.5
$wsTopo.Cells.Item(10,7).Value = 1
$wsTopo.Cells.Item(10,8).Value = 1

# ---------------------------------------------------------------------------
# Sheet "cf" (sheet6): add pv2 capacity-factor columns E/F/G (mirrors pv1's
# B/C/D: a literal value column plus two "=1*E#" computed columns).
# ---------------------------------------------------------------------------
$wsCf.Range("E1:G25").ClearFormats()
$wsCf.Cells.Item(1,5).Value = "pv2,s1"
$wsCf.Cells.Item(1,6).Value = "pv2,s2"
$wsCf.Cells.Item(1,7).Value = "pv2,s3"

for ($r = 2; $r -le 25; $r++) {
    $bVal = $wsCf.Cells.Item($r,2).Value2()
    $wsCf.Cells.Item($r,5).Value = $bVal
    $wsCf.Cells.Item($r,6).Formula = "=1*E$r"
    $wsCf.Cells.Item($r,7).Formula = "=1*E$r"
}

# ---------------------------------------------------------------------------
# View-state cosmetics: selected cell per sheet + active tab (matches the
# author ending on the "processes" sheet with D10 selected).
# ---------------------------------------------------------------------------
$wsNodes = $wb.Worksheets.Item("nodes")
$wsNodes.Range("C3").Select()

$wsTopo.Range("B15").Select()

$wsCf.Range("I8").Select()

$wsProc.Activate()
$wsProc.Range("D10").Select()
